$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new values look numeric (e.g. "1.00", "322.99").
# Coinranking prices are stored as text labels in this sheet, so we force
# Text format before writing such values to stop Excel from re-typing them
# as numbers, then restore General/Normal so no stray style is left on the
# cell. Values that can never parse as a number (e.g. "47.535.34", which has
# two decimal points) are simply assigned directly.
function Set-PriceText($cellRef, $text) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.NumberFormat = "General"
    $cell.Style = "Normal"
}

$ws.Range('D2').Value = '47.535.34'
$ws.Range('E2').Value = '  +4.83%  '

$ws.Range('D3').Value = '2.486.58'
$ws.Range('E3').Value = '  +2.51%  '

Set-PriceText 'D4' '1.00'
$ws.Range('E4').Value = '  +0.08%  '

Set-PriceText 'D5' '322.99'
$ws.Range('E5').Value = '  +1.50%  '

Set-PriceText 'D6' '105.11'
$ws.Range('E6').Value = '  +2.59%  '

Set-PriceText 'D7' '0.525'
$ws.Range('E7').Value = '  +1.75%  '

Set-PriceText 'D8' '1.00'
$ws.Range('E8').Value = '  +0.01%  '

Set-PriceText 'D9' '0.542'
$ws.Range('E9').Value = '  +2.13%  '

Set-PriceText 'D10' '38.01'
$ws.Range('E10').Value = '  +6.79%  '

Set-PriceText 'D11' '0.0813'
$ws.Range('E11').Value = '  +1.19%  '

Set-PriceText 'D12' '0.123'
$ws.Range('E12').Value = '  +1.04%  '

Set-PriceText 'D13' '18.22'
$ws.Range('E13').Value = '  +0.39%  '

Set-PriceText 'D14' '7.15'
$ws.Range('E14').Value = '  +1.62%  '

$ws.Range('D15').Value = '2.877.12'
$ws.Range('E15').Value = '  +2.61%  '

$ws.Range('D16').Value = '2.490.87'
$ws.Range('E16').Value = '  +2.53%  '

Set-PriceText 'D17' '0.842'
$ws.Range('E17').Value = '  -0.04%  '

$ws.Range('D18').Value = '47.427.18'
$ws.Range('E18').Value = '  +4.81%  '

Set-PriceText 'D19' '12.66'
$ws.Range('E19').Value = '  +3.56%  '

Set-PriceText 'D20' '6.55'
$ws.Range('E20').Value = '  +3.26%  '

$ws.Range('D21').Value = '0.0₃0935'
$ws.Range('E21').Value = '  +1.65%  '

Set-PriceText 'D22' '70.67'
$ws.Range('E22').Value = '  +2.66%  '

Set-PriceText 'D23' '250.84'
$ws.Range('E23').Value = '  +2.73%  '

Set-PriceText 'D24' '2.38'
$ws.Range('E24').Value = '  +5.68%  '

$ws.Range('E25').Value = '  +2.70%  '

Set-PriceText 'D26' '26.11'
$ws.Range('E26').Value = '  +2.11%  '

$ws.Range('E27').Value = '  -0.09%  '

Set-PriceText 'D28' '10.00'
$ws.Range('E28').Value = '  +4.30%  '

$ws.Range('E29').Value = '  +6.59%  '

Set-PriceText 'D30' '34.94'
$ws.Range('E30').Value = '  +6.23%  '

$ws.Range('E31').Value = '  +6.07%  '

Set-PriceText 'D32' '49.48'
$ws.Range('E32').Value = '  +0.79%  '

Set-PriceText 'D33' '19.83'
$ws.Range('E33').Value = '  -1.79%  '

Set-PriceText 'D34' '5.34'
$ws.Range('E34').Value = '  +2.56%  '

$ws.Range('E35').Value = '  +2.05%  '

$ws.Range('E36').Value = '  +0.21%  '

$ws.Range('E37').Value = '  +3.46%  '

Set-PriceText 'D38' '4.59'
$ws.Range('E38').Value = '  +2.97%  '

$ws.Range('E39').Value = '  +4.11%  '

$ws.Range('E40').Value = '  +1.81%  '

$ws.Range('E41').Value = '  +1.65%  '

Set-PriceText 'D42' '121.00'
$ws.Range('E42').Value = '  -4.65%  '

Set-PriceText 'D43' '21.18'
$ws.Range('E43').Value = '  +2.63%  '

$ws.Range('E44').Value = '  +2.64%  '

$ws.Range('D45').Value = '1.961.52'
$ws.Range('E45').Value = '  +1.39%  '

Set-PriceText 'D46' '2.96'
$ws.Range('E46').Value = '  +1.26%  '

$ws.Range('E48').Value = '  +1.18%  '

$ws.Range('E49').Value = '  -1.59%  '

Set-PriceText 'D50' '5.30'
$ws.Range('E50').Value = '  +12.19%  '

Set-PriceText 'D51' '79.40'
$ws.Range('E51').Value = '  +3.87%  '
